$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Oceania" data block currently sitting in row 7 needs to move down so
# it ends up after the North America rows (i.e. becomes row 11), with the
# Europe / North America rows that followed it shifting up one row to fill
# the gap (rows 8-11 -> rows 7-10).

# Save the Oceania row (A7 and C7:J7 - B is just a spacer/style column) before
# it gets overwritten by the shift.
$savedA = $ws.Range("A7").Value2
$savedRow = $ws.Range("C7:J7").Value2

# Shift rows 8-11 up into rows 7-10.
for ($r = 8; $r -le 11; $r++) {
    $srcA = $ws.Range("A" + $r).Value2
    $srcRow = $ws.Range("C" + $r + ":J" + $r).Value2
    $dest = $r - 1
    $ws.Range("A" + $dest).Value2 = $srcA
    $ws.Range("C" + $dest + ":J" + $dest).Value2 = $srcRow
}

# Write the saved Oceania row into its new home, row 11.
$ws.Range("A11").Value2 = $savedA
$ws.Range("C11:J11").Value2 = $savedRow

# Update the active selection to match the authored state.
$ws.Range("D21").Select()
